$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 585 by copying it and inserting the copy above itself,
# shifting rows 585-624 down to 586-625 (new row 625 created).
$ws.Rows("585:585").Copy() | Out-Null
$ws.Rows("585:585").Insert() | Out-Null
